$wb = $excel.ActiveWorkbook

# Overview sheet: mark the 8ccf... file as handed back (same status as the 85de... row)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: update status + handback datetime for 8ccf... row, and refresh handback datetime for 85de... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-03-20 17:40:56"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-20 17:40:56"

# de-de sheet: update status + handback datetime for 8ccf... row, and refresh handback datetime for 85de... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-03-20 17:41:10"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-20 17:41:10"
